$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
